$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates per diff
$ws.Range("N2").Value = "2019-12-31 00:00:00"
$ws.Range("O2").Value = 285454044.81
$ws.Range("P2").Value = 4287579921.77
$ws.Range("Q2").Value = 3978901676.27
$ws.Range("R2").Value = 8.2476503683
$ws.Range("S2").Value = 1675629906.3
$ws.Range("T2").Value = 1675629906.3
$ws.Range("U2").Value = 5.6774479736
$ws.Range("V2").Value = 1881247724.37
$ws.Range("W2").Value = 181051273.92
$ws.Range("X2").Value = 38165594.12
$ws.Range("Y2").Value = 328058787.31
$ws.Range("Z2").Value = 331767472.85
$ws.Range("AA2").Value = 43973186.77
$ws.Range("AG2").Value = 51814533.79
$ws.Range("AP2").Value = 7.6550376209
$ws.Range("AQ2").Value = 1.437314814997
$ws.Range("AR2").Value = 6.532937926962
$ws.Range("AS2").Value = 268811752.24
$ws.Range("AT2").Value = 3.13044955459
